# Generate Report for Handoff
#
# The localization-status report was regenerated: the "602bb916..." file
# advanced past the "b4e0931c..." file in the handoff queue, so their rows
# swap places (and b4e0931c's status flips from "In Translation" to
# "Ready for handoff"). The "Latest Handoff Datetime" for the zh-cn / de-de
# sheets is bumped to the new handoff run's timestamp everywhere except the
# still-untouched "93843a23..." row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A5").Value2 = "602bb916-6ec5-48b6-9613-1b1b4cc0111b.md"
$ws.Range("B5").Value2 = "Ready for handoff"
$ws.Range("C5").Value2 = "Ready for handoff"
$ws.Range("A7").Value2 = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.md"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$5') {
        $h.TextToDisplay = "602bb916-6ec5-48b6-9613-1b1b4cc0111b.md"
    }
    elseif ($addr -eq '$A$7') {
        $h.TextToDisplay = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("D2").Value2 = "2016-03-03 13:29:10"
$ws.Range("D3").Value2 = "2016-03-03 13:29:10"

$ws.Range("A5").Value2 = "602bb916-6ec5-48b6-9613-1b1b4cc0111b.md"
$ws.Range("B5").Value2 = "Ready for handoff"
$ws.Range("C5").Value2 = "602bb916-6ec5-48b6-9613-1b1b4cc0111b.e1fd8e08508ba4aeaac17ed8af323330436f6dca.zh-cn.xlf"
$ws.Range("D5").Value2 = "2016-03-03 13:29:10"

$ws.Range("D6").Value2 = "2016-03-03 13:29:10"

$ws.Range("A7").Value2 = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.md"
$ws.Range("C7").Value2 = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.8261d65499c1c71b25bcb0f4c39c0a00967bb7b1.zh-cn.xlf"
$ws.Range("D7").Value2 = "2016-03-03 13:29:10"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$5') {
        $h.TextToDisplay = "602bb916-6ec5-48b6-9613-1b1b4cc0111b.md"
    }
    elseif ($addr -eq '$C$5') {
        $h.TextToDisplay = "602bb916-6ec5-48b6-9613-1b1b4cc0111b.e1fd8e08508ba4aeaac17ed8af323330436f6dca.zh-cn.xlf"
    }
    elseif ($addr -eq '$A$7') {
        $h.TextToDisplay = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.md"
    }
    elseif ($addr -eq '$C$7') {
        $h.TextToDisplay = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.8261d65499c1c71b25bcb0f4c39c0a00967bb7b1.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("D2").Value2 = "2016-03-03 13:29:27"
$ws.Range("D3").Value2 = "2016-03-03 13:29:27"

$ws.Range("A5").Value2 = "602bb916-6ec5-48b6-9613-1b1b4cc0111b.md"
$ws.Range("B5").Value2 = "Ready for handoff"
$ws.Range("C5").Value2 = "602bb916-6ec5-48b6-9613-1b1b4cc0111b.e1fd8e08508ba4aeaac17ed8af323330436f6dca.de-de.xlf"
$ws.Range("D5").Value2 = "2016-03-03 13:29:27"

$ws.Range("D6").Value2 = "2016-03-03 13:29:27"

$ws.Range("A7").Value2 = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.md"
$ws.Range("C7").Value2 = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.8261d65499c1c71b25bcb0f4c39c0a00967bb7b1.de-de.xlf"
$ws.Range("D7").Value2 = "2016-03-03 13:29:27"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$5') {
        $h.TextToDisplay = "602bb916-6ec5-48b6-9613-1b1b4cc0111b.md"
    }
    elseif ($addr -eq '$C$5') {
        $h.TextToDisplay = "602bb916-6ec5-48b6-9613-1b1b4cc0111b.e1fd8e08508ba4aeaac17ed8af323330436f6dca.de-de.xlf"
    }
    elseif ($addr -eq '$A$7') {
        $h.TextToDisplay = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.md"
    }
    elseif ($addr -eq '$C$7') {
        $h.TextToDisplay = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.8261d65499c1c71b25bcb0f4c39c0a00967bb7b1.de-de.xlf"
    }
}
